$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.352"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05631"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.430"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.372"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8177"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9204"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1445"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07441"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03249"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03090"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09314"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001639"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04737"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006406"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005063"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001032"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.166"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01147"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("D25").Value = "'0.3306"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1319"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "26AAXTokenAABWorstin24h"
$ws.Range("D28").Value = "'0.0002999"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03939"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1066"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003400"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.008521"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005572"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Value = "'0.6798"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.1888"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.01010"
$ws.Range("D51").Style = "Normal"
